# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-199) from serial date 45188 (2023-09-19) to 45189 (2023-09-20).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C199").Value = 45189
